{"js": "// Fixed #542 Add template visibility.\n// The `{m:template myTemplate(...)}` field is changed to\n// `{m:template public myTemplate(...)}` by inserting the \"public\"\n// visibility keyword between the \"m:template\" keyword and the\n// template name \"myTemplate\".\n\nconst searchResults = context.document.body.search(\"m:template myTemplate\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the 'm:template myTemplate' field text to update.\");\n}\n\n// There is a single occurrence of this field in the document; replace its\n// text in place (this preserves any surrounding structure, e.g. bookmarks,\n// that sit inside the same paragraph but outside the matched range).\nsearchResults.items[0].insertText(\"m:template public myTemplate\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Fixed #542 Add template visibility.\n# The `{m:template myTemplate(...)}` field is changed to\n# `{m:template public myTemplate(...)}` by inserting the \"public\"\n# visibility keyword between the \"m:template\" keyword and the\n# template name \"myTemplate\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"m:template myTemplate\",  # FindText\n    $false,                   # MatchCase\n    $false,                   # MatchWholeWord\n    $false,                   # MatchWildcards\n    $false,                   # MatchSoundsLike\n    $false,                   # MatchAllWordForms\n    $true,                    # Forward\n    1,                        # Wrap (wdFindContinue)\n    $false,                   # Format\n    \"m:template public myTemplate\",  # ReplaceWith\n    2                         # Replace (wdReplaceAll)\n)\n"}
